$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update a few values on row 4 ---
$ws.Range("C4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0.3
$ws.Range("AC4").Value = 0.25
$ws.Range("AD4").Value = 0.4285714285714285
$ws.Range("AE4").Value = 0.8
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = 4

# --- Insert a new row 48 (student 211062992) shifting existing rows 48-52 down to 49-53 ---
$ws.Rows.Item(48).Insert()

# Match the look (bold, centered, bordered) of the other ID cells in column A,
# then force text formatting so the ID is stored as a string like its siblings.
$ws.Range("A47").Copy()
$ws.Range("A48").PasteSpecial(-4122)
$ws.Range("A48").NumberFormat = "@"
$ws.Range("A48").Value = "211062992"
$newRowValues = New-Object 'object[,]' 1,32
for ($i = 0; $i -lt 32; $i++) { $newRowValues[0, $i] = 0 }
$ws.Range("B48:AG48").Value = $newRowValues
$ws.Range("AH48").Value = "II"

# --- Row 53 (previously row 52, shifted down by the insert) gets a few updated values ---
$ws.Range("L53").Value = 2
$ws.Range("M53").Value = 1
$ws.Range("N53").Value = 1
$ws.Range("AG53").Value = 25
$ws.Range("AH53").Value = "SS"
